# 4.0.3 model and data
# Updates the "Boolean" sheet's list of constrained-data-type InputData CSVs:
#   - "trans/BVTQaZ/BVTQaZ.csv" is split into six mode-specific files
#   - "trans/VTQaZ/VTQaZ.csv" is split into six mode-specific files
# and restores "About" as the active/selected sheet.

$wb = $excel.ActiveWorkbook

$wsAbout    = $wb.Worksheets.Item("About")
$wsInteger  = $wb.Worksheets.Item("Integer")
$wsBoolean  = $wb.Worksheets.Item("Boolean")
$wsSubscript = $wb.Worksheets.Item("Subscript")

# ---------------------------------------------------------------------------
# Boolean sheet: expand the two aggregate "trans" rows into per-mode rows.
# ---------------------------------------------------------------------------

# Row 21 (trans/VTQaZ/VTQaZ.csv) is split first so row numbers for the
# first split (row 17) aren't disturbed by this one.
$wsBoolean.Rows("21:21").Insert()
$wsBoolean.Rows("21:21").Insert()
$wsBoolean.Rows("21:21").Insert()
$wsBoolean.Rows("21:21").Insert()
$wsBoolean.Rows("21:21").Insert()

$wsBoolean.Range("A21").Value = "trans/VTQaZ/VTQaZ-LDVs.csv"
$wsBoolean.Range("A22").Value = "trans/VTQaZ/VTQaZ-HDVs.csv"
$wsBoolean.Range("A23").Value = "trans/VTQaZ/VTQaZ-aircraft.csv"
$wsBoolean.Range("A24").Value = "trans/VTQaZ/VTQaZ-rail.csv"
$wsBoolean.Range("A25").Value = "trans/VTQaZ/VTQaZ-ships.csv"
$wsBoolean.Range("A26").Value = "trans/VTQaZ/VTQaZ-motorbikes.csv"

# Row 17 (trans/BVTQaZ/BVTQaZ.csv)
$wsBoolean.Rows("17:17").Insert()
$wsBoolean.Rows("17:17").Insert()
$wsBoolean.Rows("17:17").Insert()
$wsBoolean.Rows("17:17").Insert()
$wsBoolean.Rows("17:17").Insert()

$wsBoolean.Range("A17").Value = "trans/BVTQaZ/BVTQaZ-LDVs.csv"
$wsBoolean.Range("A18").Value = "trans/BVTQaZ/BVTQaZ-HDVs.csv"
$wsBoolean.Range("A19").Value = "trans/BVTQaZ/BVTQaZ-aircraft.csv"
$wsBoolean.Range("A20").Value = "trans/BVTQaZ/BVTQaZ-rail.csv"
$wsBoolean.Range("A21").Value = "trans/BVTQaZ/BVTQaZ-ships.csv"
$wsBoolean.Range("A22").Value = "trans/BVTQaZ/BVTQaZ-motorbikes.csv"

# Six trailing blank (but formatted) rows after the data, as in the source.
$wsBoolean.Rows("33:38").Font.Name = $wsBoolean.Range("A2").Font.Name

# ---------------------------------------------------------------------------
# View-state: make "About" the active tab again, and restore the cursor
# positions left on "Integer" and "Boolean" by the edits above.
# ---------------------------------------------------------------------------

$wsInteger.Range("A13").Select()
$wsBoolean.Range("A32").Select()
$wsAbout.Activate()
$wsAbout.Range("A1").Select()

Write-Host "Boolean sheet InputData list updated; About is active tab."
